# Update countries & provincias Spain
# - Swap the display order of 4 neighbouring country rows (their case
#   counts effectively changed rank order vs. each other).
# - Refresh the "datos actualizados" timestamp.
# - Refresh the latest case/recovered/death counters for the countries
#   whose figures moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap adjacent country names (rank reorder) ------------------------
# Armenia / Nigeria (rows 52-53)
$ws.Range("A52").Value = "Nigeria"
$ws.Range("A53").Value = "Armenia"

# Nepal / Chequia (rows 68-69)
$ws.Range("A68").Value = "Chequia"
$ws.Range("A69").Value = "Nepal"

# Dominica / Fiyi (rows 202-203)
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"

# Groenlandia / Islas Malvinas (rows 208-209)
$ws.Range("A208").Value = "Islas Malvinas"
$ws.Range("A209").Value = "Groenlandia"

# --- 2) Refresh the "updated at" banner ------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Junio de 2020 a las 01:21"

# --- 3) Refresh updated statistics -----------------------------------------
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2461439
$ws.Range("C4").Value = 37271
$ws.Range("D4").Value = 1033914
$ws.Range("E4").Value = 1303278
$ws.Range("G4").Value = 774
$ws.Range("H4").Value = 124247

# Row 5 - Brasil
$ws.Range("B5").Value = 1192474
$ws.Range("C5").Value = 40995
$ws.Range("E5").Value = 488692
$ws.Range("G5").Value = 1103
$ws.Range("H5").Value = 53874

# Row 25 - Colombia
$ws.Range("B25").Value = 77113
$ws.Range("C25").Value = 3541
$ws.Range("D25").Value = 31671
$ws.Range("E25").Value = 42951
$ws.Range("G25").Value = 87
$ws.Range("H25").Value = 2491

# Row 52 - now Nigeria
$ws.Range("B52").Value = 22020
$ws.Range("C52").Value = 649
$ws.Range("D52").Value = 7613
$ws.Range("E52").Value = 13865
$ws.Range("G52").Value = 9
$ws.Range("H52").Value = 542

# Row 53 - now Armenia
$ws.Range("B53").Value = 21717
$ws.Range("C53").Value = 711
$ws.Range("D53").Value = 10797
$ws.Range("E53").Value = 10534
$ws.Range("G53").Value = 14
$ws.Range("H53").Value = 386

# Row 55 - Japon
$ws.Range("B55").Value = 18024
$ws.Range("C55").Value = 56
$ws.Range("D55").Value = 16263
$ws.Range("E55").Value = 798
$ws.Range("G55").Value = 8
$ws.Range("H55").Value = 963

# Row 68 - now Chequia
$ws.Range("B68").Value = 10777
$ws.Range("C68").Value = 127
$ws.Range("D68").Value = 7588
$ws.Range("E68").Value = 2846
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 343

# Row 69 - now Nepal
$ws.Range("C69").Value = 629
$ws.Range("D69").Value = 2338
$ws.Range("E69").Value = 8366
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 24

# Row 86 - Gabon
$ws.Range("B86").Value = 4956
$ws.Range("C86").Value = 107
$ws.Range("D86").Value = 2177
$ws.Range("E86").Value = 2740

# Row 132 - Republica de Chipre
$ws.Range("B132").Value = 991
$ws.Range("C132").Value = 1
$ws.Range("E132").Value = 148

# Row 137 - Uruguay
$ws.Range("B137").Value = 902
$ws.Range("C137").Value = 17
$ws.Range("E137").Value = 61
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 26

# Row 151 - Togo
$ws.Range("B151").Value = 583
$ws.Range("C151").Value = 7
$ws.Range("D151").Value = 392
$ws.Range("E151").Value = 177
$ws.Range("G151").Value = 1
$ws.Range("H151").Value = 14

# Row 166 - Guyana
$ws.Range("B166").Value = 209
$ws.Range("C166").Value = 3
$ws.Range("E166").Value = 90

# Row 183 - Liechtenstein
$ws.Range("D183").Value = 81
$ws.Range("E183").Value = 0
